$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2,1).Value = "Última actualización: 04:48:57"
$ws.Cells.Item(3,1).Value = "Total filas: 34"
$ws.Cells.Item(17,1).Value = "04:48:57"
$ws.Cells.Item(17,4).Value = 5
$ws.Cells.Item(18,1).Value = "04:48:57"
$ws.Cells.Item(18,2).Value = "05:14"
$ws.Cells.Item(18,4).Value = 26
$ws.Cells.Item(19,2).Value = "05:15"
$ws.Cells.Item(19,3).Value = "14_ABASTO"
$ws.Cells.Item(19,4).Value = 40
$ws.Cells.Item(20,1).Value = "04:48:57"
$ws.Cells.Item(20,2).Value = "05:16"
$ws.Cells.Item(20,3).Value = "17_ROMERO"
$ws.Cells.Item(20,4).Value = 28
$ws.Cells.Item(21,1).Value = "04:48:57"
$ws.Cells.Item(21,2).Value = "05:21"
$ws.Cells.Item(21,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(21,4).Value = 33
$ws.Cells.Item(22,2).Value = "05:22"
$ws.Cells.Item(22,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(22,4).Value = 47
$ws.Cells.Item(23,1).Value = "03:35:49"
$ws.Cells.Item(23,2).Value = "05:28"
$ws.Cells.Item(23,3).Value = "14_ABASTO"
$ws.Cells.Item(23,4).Value = 113
$ws.Cells.Item(24,1).Value = "04:48:57"
$ws.Cells.Item(24,2).Value = "05:34"
$ws.Cells.Item(24,3).Value = "215B_EL PATO"
$ws.Cells.Item(24,4).Value = 46
$ws.Cells.Item(25,1).Value = "04:01:13"
$ws.Cells.Item(25,2).Value = "05:35"
$ws.Cells.Item(25,3).Value = "215B_EL PATO"
$ws.Cells.Item(25,4).Value = 94
$ws.Cells.Item(26,1).Value = "04:01:13"
$ws.Cells.Item(26,2).Value = "05:37"
$ws.Cells.Item(26,3).Value = "14_ABASTO"
$ws.Cells.Item(26,4).Value = 96
$ws.Cells.Item(27,1).Value = "04:48:57"
$ws.Cells.Item(27,2).Value = "05:46"
$ws.Cells.Item(27,3).Value = "15_ABASTO"
$ws.Cells.Item(27,4).Value = 58
$ws.Cells.Item(28,1).Value = "04:48:57"
$ws.Cells.Item(28,2).Value = "06:04"
$ws.Cells.Item(28,3).Value = "16_SANTA ANA"
$ws.Cells.Item(28,4).Value = 76
$ws.Cells.Item(29,1).Value = "04:48:57"
$ws.Cells.Item(29,2).Value = "06:11"
$ws.Cells.Item(29,3).Value = "215A_EL PATO"
$ws.Cells.Item(29,4).Value = 83
$ws.Cells.Item(30,1).Value = "04:48:57"
$ws.Cells.Item(30,2).Value = "06:13"
$ws.Cells.Item(30,3).Value = "225_HARAS DEL SUR"
$ws.Cells.Item(30,4).Value = 85
$ws.Cells.Item(31,2).Value = "06:14"
$ws.Cells.Item(31,3).Value = "225_HARAS DEL SUR"
$ws.Cells.Item(31,4).Value = 99
$ws.Cells.Item(32,1).Value = "04:48:57"
$ws.Cells.Item(32,2).Value = "06:20"
$ws.Cells.Item(32,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(32,4).Value = 92
$ws.Cells.Item(33,1).Value = "04:35:25"
$ws.Cells.Item(33,2).Value = "06:21"
$ws.Cells.Item(33,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(33,4).Value = 106
$ws.Cells.Item(33,5).Value = "LP1912"
$ws.Cells.Item(34,1).Value = "04:48:57"
$ws.Cells.Item(34,2).Value = "06:26"
$ws.Cells.Item(34,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(34,4).Value = 98
$ws.Cells.Item(34,5).Value = "LP1912"
$ws.Cells.Item(35,1).Value = "04:35:25"
$ws.Cells.Item(35,2).Value = "06:27"
$ws.Cells.Item(35,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(35,4).Value = 112
$ws.Cells.Item(35,5).Value = "LP1912"
$ws.Cells.Item(36,1).Value = "04:48:57"
$ws.Cells.Item(36,2).Value = "06:29"
$ws.Cells.Item(36,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(36,4).Value = 101
$ws.Cells.Item(36,5).Value = "LP1912"
$ws.Cells.Item(37,1).Value = "04:48:57"
$ws.Cells.Item(37,2).Value = "06:31"
$ws.Cells.Item(37,3).Value = "16_SANTA ANA"
$ws.Cells.Item(37,4).Value = 103
$ws.Cells.Item(37,5).Value = "LP1912"
$ws.Cells.Item(38,1).Value = "04:48:57"
$ws.Cells.Item(38,2).Value = "06:43"
$ws.Cells.Item(38,3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(38,4).Value = 115
$ws.Cells.Item(38,5).Value = "LP1912"
$ws.Cells.Item(39,1).Value = "04:48:57"
$ws.Cells.Item(39,2).Value = "06:46"
$ws.Cells.Item(39,3).Value = "215C_EL PATO"
$ws.Cells.Item(39,4).Value = 118
$ws.Cells.Item(39,5).Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2,1).Value = "Última actualización: 04:48:57"
$ws.Cells.Item(3,1).Value = "Total filas: 10"
$ws.Cells.Item(12,1).Value = "04:48:57"
$ws.Cells.Item(12,4).Value = 46
$ws.Cells.Item(14,1).Value = "04:48:57"
$ws.Cells.Item(14,4).Value = 83
$ws.Cells.Item(15,1).Value = "04:48:57"
$ws.Cells.Item(15,2).Value = "06:46"
$ws.Cells.Item(15,3).Value = "215C_EL PATO"
$ws.Cells.Item(15,4).Value = 118
$ws.Cells.Item(15,5).Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2,1).Value = "Última actualización: 04:48:57"
$ws.Cells.Item(3,1).Value = "Total filas: 7"
$ws.Cells.Item(7,1).Value = "04:48:57"
$ws.Cells.Item(7,2).Value = "05:43"
$ws.Cells.Item(7,4).Value = 55
$ws.Cells.Item(8,2).Value = "05:44"
$ws.Cells.Item(8,4).Value = 69
$ws.Cells.Item(9,1).Value = "04:48:57"
$ws.Cells.Item(9,2).Value = "06:08"
$ws.Cells.Item(9,3).Value = "215A_LA PLATA"
$ws.Cells.Item(9,4).Value = 80
$ws.Cells.Item(9,5).Value = "L6173"
$ws.Cells.Item(10,1).Value = "04:35:25"
$ws.Cells.Item(10,2).Value = "06:09"
$ws.Cells.Item(10,3).Value = "215A_LA PLATA"
$ws.Cells.Item(10,4).Value = 94
$ws.Cells.Item(10,5).Value = "L6173"
$ws.Cells.Item(11,1).Value = "04:48:57"
$ws.Cells.Item(11,2).Value = "06:32"
$ws.Cells.Item(11,3).Value = "215C_LA PLATA"
$ws.Cells.Item(11,4).Value = 104
$ws.Cells.Item(11,5).Value = "L6203"
$ws.Cells.Item(12,1).Value = "04:35:25"
$ws.Cells.Item(12,2).Value = "06:33"
$ws.Cells.Item(12,3).Value = "215C_LA PLATA"
$ws.Cells.Item(12,4).Value = 118
$ws.Cells.Item(12,5).Value = "L6203"
